$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared string text used by E26 (remove trailing sentence about loading UI)
$ws.Range("E26").Value = "Redid UI to more closely match Excel ribbon idioms. Used bootstrap for element styles."

# Update selection to H27
$ws.Range("H27").Select()

# Row 26: reset row height to default (remove explicit 25.5 height) and update status/date
$ws.Rows("26").AutoFit()

# F26: change status from "In progress" to "Complete"
$ws.Range("F26").Value = "Complete"

# H26: set completed date (use same formatting as G26, a date-formatted cell)
$ws.Range("G26").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("H26").Value2 = 42984
$excel.CutCopyMode = $false
